$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.032.60"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "'1.621.62"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'215.08"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").Value = "'20.16"
$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'1.636.86"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").Value = "'4.12"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("E15").Value = "  -4.68%  "

$ws.Range("D16").Value = "'27.015.93"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "'0.0₃0737"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "'215.78"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").Value = "'6.89"
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("E22").Value = "  -6.17%  "

$ws.Range("D23").Value = "'8.97"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").Value = "'147.45"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "'7.27"
$ws.Range("E26").Value = "  -4.04%  "

$ws.Range("D27").Value = "'0.118"
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("D28").Value = "'15.55"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D33").Value = "'1.335.09"
$ws.Range("E33").Value = "  +5.72%  "

$ws.Range("E34").Value = "  -1.34%  "

$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("D36").Value = "'0.0176"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("D37").Value = "'0.544"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("D40").Value = "'2.24"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "'0.800"
$ws.Range("E41").Value = "  -0.98%  "

$ws.Range("D42").Value = "'64.30"
$ws.Range("E42").Value = "  +3.93%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.22"
$ws.Range("E43").Value = "  -3.96%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "'1.758.66"
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("D45").Value = "'90.36"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("D48").Value = "'0.820"
$ws.Range("E48").Value = "  +22.28%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").Value = "'7.52"
$ws.Range("E51").Value = "  -1.74%  "
